$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 326387
$ws.Range("C2").Value = 292481
$ws.Range("D2").Value = 314108
$ws.Range("E2").Value = 313198
$ws.Range("F2").Value = 307805
$ws.Range("G2").Value = 284772
$ws.Range("H2").Value = 302068
$ws.Range("I2").Value = 295446
$ws.Range("J2").Value = 292653
$ws.Range("K2").Value = 316233
$ws.Range("L2").Value = 302574
$ws.Range("M2").Value = 332356
$ws.Range("N2").Value = 3680081

$ws.Range("B3").Value = 24343
$ws.Range("C3").Value = 21038
$ws.Range("D3").Value = 24972
$ws.Range("E3").Value = 27290
$ws.Range("F3").Value = 26278
$ws.Range("G3").Value = 14652
$ws.Range("H3").Value = 19440
$ws.Range("I3").Value = 20496
$ws.Range("J3").Value = 16979
$ws.Range("K3").Value = 15918
$ws.Range("L3").Value = 11188
$ws.Range("M3").Value = 11945
$ws.Range("N3").Value = 234539

$ws.Range("B4").Value = 36532
$ws.Range("C4").Value = 36992
$ws.Range("D4").Value = 38227
$ws.Range("E4").Value = 40459
$ws.Range("F4").Value = 40148
$ws.Range("G4").Value = 28873
$ws.Range("H4").Value = 36623
$ws.Range("I4").Value = 38352
$ws.Range("J4").Value = 36974
$ws.Range("K4").Value = 42030
$ws.Range("L4").Value = 38401
$ws.Range("M4").Value = 34125
$ws.Range("N4").Value = 447736

$ws.Range("B5").Value = 17128
$ws.Range("C5").Value = 17363
$ws.Range("D5").Value = 19584
$ws.Range("E5").Value = 18609
$ws.Range("F5").Value = 18258
$ws.Range("G5").Value = 22252
$ws.Range("H5").Value = 19119
$ws.Range("I5").Value = 19711
$ws.Range("J5").Value = 26336
$ws.Range("K5").Value = 23936
$ws.Range("L5").Value = 24282
$ws.Range("M5").Value = 28839
$ws.Range("N5").Value = 255417

$ws.Range("B6").Value = 7217
$ws.Range("C6").Value = 19335
$ws.Range("D6").Value = 17412
$ws.Range("E6").Value = 13491
$ws.Range("F6").Value = 19251
$ws.Range("G6").Value = 16131
$ws.Range("H6").Value = 20046
$ws.Range("I6").Value = 17582
$ws.Range("J6").Value = 17344
$ws.Range("K6").Value = 20299
$ws.Range("L6").Value = 16931
$ws.Range("M6").Value = 15946
$ws.Range("N6").Value = 200985

$ws.Range("B7").Value = 15317
$ws.Range("C7").Value = 14421
$ws.Range("D7").Value = 14850
$ws.Range("E7").Value = 18373
$ws.Range("F7").Value = 16362
$ws.Range("G7").Value = 16696
$ws.Range("H7").Value = 15999
$ws.Range("I7").Value = 16723
$ws.Range("J7").Value = 14182
$ws.Range("K7").Value = 16410
$ws.Range("L7").Value = 16483
$ws.Range("M7").Value = 13552
$ws.Range("N7").Value = 189368

$ws.Range("B8").Value = 6970
$ws.Range("C8").Value = 7003
$ws.Range("D8").Value = 11414
$ws.Range("E8").Value = 16388
$ws.Range("F8").Value = 18673
$ws.Range("G8").Value = 20900
$ws.Range("H8").Value = 21377
$ws.Range("I8").Value = 21668
$ws.Range("J8").Value = 19402
$ws.Range("K8").Value = 14215
$ws.Range("L8").Value = 24054
$ws.Range("M8").Value = 36291
$ws.Range("N8").Value = 218355

$ws.Range("B9").Value = 6699
$ws.Range("C9").Value = 6260
$ws.Range("D9").Value = 7951
$ws.Range("E9").Value = 7762
$ws.Range("F9").Value = 6688
$ws.Range("G9").Value = 6517
$ws.Range("H9").Value = 7616
$ws.Range("I9").Value = 8714
$ws.Range("J9").Value = 7783
$ws.Range("K9").Value = 5716
$ws.Range("L9").Value = 5096
$ws.Range("M9").Value = 5061
$ws.Range("N9").Value = 81863

$ws.Range("B10").Value = 7441
$ws.Range("C10").Value = 6524
$ws.Range("D10").Value = 6724
$ws.Range("E10").Value = 6483
$ws.Range("F10").Value = 6391
$ws.Range("G10").Value = 5799
$ws.Range("H10").Value = 6511
$ws.Range("I10").Value = 6415
$ws.Range("J10").Value = 6436
$ws.Range("K10").Value = 6877
$ws.Range("L10").Value = 6791
$ws.Range("M10").Value = 8062
$ws.Range("N10").Value = 80454

$ws.Range("B11").Value = 5404
$ws.Range("C11").Value = 5010
$ws.Range("D11").Value = 5548
$ws.Range("E11").Value = 4693
$ws.Range("F11").Value = 5348
$ws.Range("G11").Value = 4415
$ws.Range("H11").Value = 6010
$ws.Range("I11").Value = 5281
$ws.Range("J11").Value = 5897
$ws.Range("K11").Value = 5850
$ws.Range("L11").Value = 5963
$ws.Range("M11").Value = 5906
$ws.Range("N11").Value = 65325

$ws.Range("B12").Value = 3634
$ws.Range("C12").Value = 3358
$ws.Range("D12").Value = 3782
$ws.Range("E12").Value = 3943
$ws.Range("F12").Value = 3787
$ws.Range("G12").Value = 2831
$ws.Range("H12").Value = 3261
$ws.Range("I12").Value = 3546
$ws.Range("J12").Value = 2886
$ws.Range("K12").Value = 2656
$ws.Range("L12").Value = 2253
$ws.Range("M12").Value = 2027
$ws.Range("N12").Value = 37964

$ws.Range("B13").Value = 3979
$ws.Range("C13").Value = 4103
$ws.Range("D13").Value = 4300
$ws.Range("E13").Value = 3973
$ws.Range("F13").Value = 4391
$ws.Range("G13").Value = 3812
$ws.Range("H13").Value = 3749
$ws.Range("I13").Value = 3872
$ws.Range("J13").Value = 3924
$ws.Range("K13").Value = 5108
$ws.Range("L13").Value = 4416
$ws.Range("M13").Value = 4143
$ws.Range("N13").Value = 49770

$ws.Range("B14").Value = 2722
$ws.Range("C14").Value = 2648
$ws.Range("D14").Value = 2792
$ws.Range("E14").Value = 2886
$ws.Range("F14").Value = 2973
$ws.Range("G14").Value = 2758
$ws.Range("H14").Value = 2861
$ws.Range("I14").Value = 2861
$ws.Range("J14").Value = 2751
$ws.Range("K14").Value = 2907
$ws.Range("L14").Value = 2710
$ws.Range("M14").Value = 2649
$ws.Range("N14").Value = 33518

$ws.Range("B15").Value = 1232
$ws.Range("C15").Value = 982
$ws.Range("D15").Value = 978
$ws.Range("E15").Value = 1022
$ws.Range("F15").Value = 1081
$ws.Range("G15").Value = 865
$ws.Range("H15").Value = 1056
$ws.Range("I15").Value = 1060
$ws.Range("J15").Value = 1022
$ws.Range("K15").Value = 1083
$ws.Range("L15").Value = 1161
$ws.Range("M15").Value = 1442
$ws.Range("N15").Value = 12984

$ws.Range("B16").Value = 995
$ws.Range("C16").Value = 581
$ws.Range("D16").Value = 1093
$ws.Range("E16").Value = 1218
$ws.Range("F16").Value = 947
$ws.Range("G16").Value = 1181
$ws.Range("H16").Value = 2514
$ws.Range("I16").Value = 2906
$ws.Range("J16").Value = 1982
$ws.Range("K16").Value = 1695
$ws.Range("L16").Value = 1432
$ws.Range("M16").Value = 1263
$ws.Range("N16").Value = 17807

$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 52
$ws.Range("I22").Value = 2850
$ws.Range("J22").Value = 2007
$ws.Range("K22").Value = 2352
$ws.Range("L22").Value = 2287
$ws.Range("M22").Value = 2544
$ws.Range("N22").Value = 12092
